$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5209010662194005
$ws.Range("C2").Value = 0.5284173549292553
$ws.Range("D2").Value = 0.3633847247877238
$ws.Range("E2").Value = 0.6028140051356835
$ws.Range("F2").Value = 0.3140407230176452
$ws.Range("G2").Value = 15
$ws.Range("B3").Value = 0.3491358409449706
$ws.Range("C3").Value = 0.3876299128740889
$ws.Range("D3").Value = 0.2064515430844597
$ws.Range("E3").Value = 0.4543693905672561
$ws.Range("F3").Value = 0.3017614833873206
$ws.Range("G3").Value = 14
$ws.Range("B4").Value = 0.2766139601088831
$ws.Range("C4").Value = 0.3205523163110667
$ws.Range("D4").Value = 0.1530262537301577
$ws.Range("E4").Value = 0.3911857023590685
$ws.Range("F4").Value = 0.287900825001641
$ws.Range("G4").Value = 13
$ws.Range("B5").Value = 0.4245700094424397
$ws.Range("C5").Value = 0.4393557528778729
$ws.Range("D5").Value = 0.2414648368369625
$ws.Range("E5").Value = 0.4913907170846459
$ws.Range("F5").Value = 0.2583974611168742
$ws.Range("G5").Value = 12
$ws.Range("B6").Value = 0.3793378286036473
$ws.Range("C6").Value = 0.4013254374874798
$ws.Range("D6").Value = 0.2006046709917273
$ws.Range("E6").Value = 0.4478891280124215
$ws.Range("F6").Value = 0.2497563433833
$ws.Range("G6").Value = 11
$ws.Range("B7").Value = 0.3288324600564539
$ws.Range("C7").Value = 0.3601674043196504
$ws.Range("D7").Value = 0.1589310008397884
$ws.Range("E7").Value = 0.3986615116107754
$ws.Range("F7").Value = 0.2375808962882351
$ws.Range("G7").Value = 10
$ws.Range("B8").Value = 0.3582448599611742
$ws.Range("C8").Value = 0.3844285013741874
$ws.Range("D8").Value = 0.1824260236518181
$ws.Range("E8").Value = 0.427113595723454
$ws.Range("F8").Value = 0.2466728085108265
$ws.Range("G8").Value = 9
$ws.Range("B9").Value = 0.3532650510030724
$ws.Range("C9").Value = 0.381336062820872
$ws.Range("D9").Value = 0.1830290214229363
$ws.Range("E9").Value = 0.4278189119509986
$ws.Range("F9").Value = 0.257976355866154
$ws.Range("G9").Value = 8
$ws.Range("B10").Value = 0.3173063146244254
$ws.Range("C10").Value = 0.3418763137882747
$ws.Range("D10").Value = 0.1463084296434522
$ws.Range("E10").Value = 0.3825028491965154
$ws.Range("F10").Value = 0.2307148046110977
$ws.Range("G10").Value = 7
$ws.Range("B11").Value = 0.3606037648954714
$ws.Range("C11").Value = 0.379698193710172
$ws.Range("D11").Value = 0.1799938953220499
$ws.Range("E11").Value = 0.4242568742189688
$ws.Range("F11").Value = 0.2448480836729456
$ws.Range("G11").Value = 6
